$d = $word.ActiveDocument

# Replace the author's display name "Billy Wade" -> "Trey Merkley"
$d.Content.Find.Execute("Billy Wade", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Trey Merkley", 2) | Out-Null

# Replace the GitHub username in the contact table
$d.Content.Find.Execute("github.com/billywade", $true, $false, $false, $false, $false,
                         $true, 1, $false, "github.com/treymerkley", 2) | Out-Null

# Update the professional summary paragraph wording
$d.Content.Find.Execute(
    "I’m a student at OSU’s satellite campus in Okmulgee majoring in software development and information security.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "I’m a graduate of OSU’s satellite campus in Okmulgee, having majored in software development and information security.",
    2) | Out-Null

# Update the bookmark name from "billy-wade" to "trey-merkley" while keeping the same range
$bm = $d.Bookmarks("billy-wade")
$rng = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("trey-merkley", $rng) | Out-Null
